$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume change (column E) values
# Prices are forced to Text format so strings such as "1.011" are not
# re-interpreted by Excel as numbers, matching the inlineStr cells in the source file.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.450.41"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -2.41%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.950.23"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.02%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.011"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.50%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "321.97"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.03%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.010"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4802"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -4.21%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.4092"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -3.28%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "53.77"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.66%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.08522"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -7.79%  "

$ws.Range("E11").Value = "  -4.12%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "22.44"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.82%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.971.67"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -1.79%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "7.581"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -4.26%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "6.164"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -4.51%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.013"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.56%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "90.50"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("E18").Value = "  -2.98%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06636"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.17%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "18.46"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -4.80%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.011"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.843"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.39%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "28.525.15"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.31%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.46"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -5.67%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.295"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.32%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.237.82"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "156.51"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "20.31"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.72%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.175"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -4.32%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.836"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -6.14%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "124.41"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.90%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.9864"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -5.86%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.09681"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.91%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.694"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.41%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.633"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.86%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.444"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -6.34%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "9.172"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.29%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02335"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -4.08%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.06197"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.73%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.252"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -3.89%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.6233"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -3.75%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "11.22"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.22%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.010"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.43%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.1921"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -3.57%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.334"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +3.30%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.5967"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.74%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "12.93"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.26%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.065"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -6.35%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "3.408"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.96%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.06822"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.32%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "110.87"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -2.10%  "
